# Scheduled market-data refresh: overwrite Universalis-derived price/profit
# columns (H:N) with freshly pulled values, per sheet/row, leaving all other
# cells (leve metadata, item ids, etc.) untouched.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6180.1816
$ws.Range("I43").Value = 6374.25
$ws.Range("J43").Value = 5662.6665
$ws.Range("K43").Value = 6374.25
$ws.Range("L43").Value = 5662.6665
$ws.Range("M43").Value = -6305.25
$ws.Range("N43").Value = -5800.6665
$ws.Range("H64").Value = 6278.6
$ws.Range("I64").Value = 6205.25
$ws.Range("J64").Value = 6305.273
$ws.Range("K64").Value = 6205.25
$ws.Range("L64").Value = 6305.273
$ws.Range("M64").Value = -5957.25
$ws.Range("N64").Value = -6801.273
$ws.Range("H67").Value = 6278.6
$ws.Range("I67").Value = 6205.25
$ws.Range("J67").Value = 6305.273
$ws.Range("K67").Value = 6205.25
$ws.Range("L67").Value = 6305.273
$ws.Range("M67").Value = -5347.25
$ws.Range("N67").Value = -8021.273
$ws.Range("H69").Value = 9851
$ws.Range("J69").Value = 9851
$ws.Range("L69").Value = 29553
$ws.Range("N69").Value = -31301
$ws.Range("H72").Value = 9851
$ws.Range("J72").Value = 9851
$ws.Range("L72").Value = 88659
$ws.Range("N72").Value = -97395
$ws.Range("H137").Value = 2318.963
$ws.Range("I137").Value = 1458.5333
$ws.Range("K137").Value = 4375.5999
$ws.Range("M137").Value = -1825.5999
$ws.Range("H138").Value = 2971.0894
$ws.Range("I138").Value = 2622.8948
$ws.Range("K138").Value = 7868.6844
$ws.Range("M138").Value = -2728.6844
$ws.Range("H141").Value = 4536.727
$ws.Range("I141").Value = 4419.476
$ws.Range("K141").Value = 13258.428
$ws.Range("M141").Value = -8078.428

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6322.705
$ws.Range("I32").Value = 5723.8965
$ws.Range("K32").Value = 5723.8965
$ws.Range("M32").Value = -5436.8965
$ws.Range("H61").Value = 5389058
$ws.Range("I61").Value = 6901169.5
$ws.Range("J61").Value = 1003936
$ws.Range("K61").Value = 6901169.5
$ws.Range("L61").Value = 1003936
$ws.Range("M61").Value = -6900957.5
$ws.Range("N61").Value = -1004360
$ws.Range("H74").Value = 1790
$ws.Range("I74").Value = 1034.9412
$ws.Range("K74").Value = 1034.9412
$ws.Range("M74").Value = -160.9412
$ws.Range("H77").Value = 1790
$ws.Range("I77").Value = 1034.9412
$ws.Range("K77").Value = 5174.706
$ws.Range("M77").Value = -806.7060000000001
$ws.Range("H103").Value = 92892
$ws.Range("J103").Value = 92892
$ws.Range("L103").Value = 92892
$ws.Range("N103").Value = -95236
$ws.Range("H128").Value = 69999.5
$ws.Range("J128").Value = 69999.5
$ws.Range("L128").Value = 69999.5
$ws.Range("N128").Value = -79959.5
$ws.Range("H132").Value = 2175293.2
$ws.Range("I132").Value = 1315.7561
$ws.Range("K132").Value = 3947.2683
$ws.Range("M132").Value = -1417.2683
$ws.Range("H136").Value = 5389058
$ws.Range("I136").Value = 6901169.5
$ws.Range("J136").Value = 1003936
$ws.Range("K136").Value = 20703508.5
$ws.Range("L136").Value = 3011808
$ws.Range("M136").Value = -20700958.5
$ws.Range("N136").Value = -3016908

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3334.9333
$ws.Range("I20").Value = 3396.2632
$ws.Range("J20").Value = 3229
$ws.Range("K20").Value = 3396.2632
$ws.Range("L20").Value = 3229
$ws.Range("M20").Value = -3149.2632
$ws.Range("N20").Value = -3723
$ws.Range("H105").Value = 517921.6
$ws.Range("I105").Value = 716436.4
$ws.Range("K105").Value = 716436.4
$ws.Range("M105").Value = -714689.4
$ws.Range("H134").Value = 3574459
$ws.Range("I134").Value = 3061.9524
$ws.Range("K134").Value = 9185.8572
$ws.Range("M134").Value = -6650.8572

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1433.3334
$ws.Range("I10").Value = 1150
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1150
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -1011
$ws.Range("N10").Value = -2278
$ws.Range("H11").Value = 18499.5
$ws.Range("I11").Value = 7000
$ws.Range("J11").Value = 29999
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 29999
$ws.Range("M11").Value = -6860
$ws.Range("N11").Value = -30279
$ws.Range("H22").Value = 256.66666
$ws.Range("I22").Value = 165
$ws.Range("K22").Value = 165
$ws.Range("M22").Value = 185
$ws.Range("H31").Value = 43481590
$ws.Range("I31").Value = 55558212
$ws.Range("K31").Value = 55558212
$ws.Range("M31").Value = -55557917
$ws.Range("H34").Value = 43481590
$ws.Range("I34").Value = 55558212
$ws.Range("K34").Value = 55558212
$ws.Range("M34").Value = -55558010
$ws.Range("H51").Value = 18793.334
$ws.Range("I51").Value = 18793.334
$ws.Range("K51").Value = 18793.334
$ws.Range("M51").Value = -18057.334
$ws.Range("H61").Value = 18793.334
$ws.Range("I61").Value = 18793.334
$ws.Range("K61").Value = 18793.334
$ws.Range("M61").Value = -18445.334

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 178.22223
$ws.Range("I2").Value = 81
$ws.Range("J2").Value = 299.75
$ws.Range("K2").Value = 486
$ws.Range("L2").Value = 1798.5
$ws.Range("M2").Value = -373
$ws.Range("N2").Value = -2024.5
$ws.Range("H109").Value = 5253.385
$ws.Range("I109").Value = 1496.1
$ws.Range("K109").Value = 4488.299999999999
$ws.Range("M109").Value = -3448.299999999999

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3032.5
$ws.Range("I80").Value = 2924.25
$ws.Range("J80").Value = 3249
$ws.Range("K80").Value = 2924.25
$ws.Range("L80").Value = 3249
$ws.Range("M80").Value = -1926.25
$ws.Range("N80").Value = -5245
$ws.Range("H83").Value = 3032.5
$ws.Range("I83").Value = 2924.25
$ws.Range("J83").Value = 3249
$ws.Range("K83").Value = 14621.25
$ws.Range("L83").Value = 16245
$ws.Range("M83").Value = -9629.25
$ws.Range("N83").Value = -26229
$ws.Range("H132").Value = 1669063.5
$ws.Range("I132").Value = 2443.8235
$ws.Range("K132").Value = 7331.470499999999
$ws.Range("M132").Value = -4801.470499999999

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17174.625
$ws.Range("J22").Value = 4316
$ws.Range("L22").Value = 4316
$ws.Range("N22").Value = -4906
$ws.Range("H27").Value = 17174.625
$ws.Range("J27").Value = 4316
$ws.Range("L27").Value = 4316
$ws.Range("N27").Value = -4530
$ws.Range("H64").Value = 60499.25
$ws.Range("J64").Value = 60499.25
$ws.Range("L64").Value = 60499.25
$ws.Range("N64").Value = -60949.25
$ws.Range("H67").Value = 60499.25
$ws.Range("J67").Value = 60499.25
$ws.Range("L67").Value = 60499.25
$ws.Range("N67").Value = -62059.25
$ws.Range("H82").Value = 4787.125
$ws.Range("I82").Value = 1059.8
$ws.Range("K82").Value = 1059.8
$ws.Range("M82").Value = -698.8
$ws.Range("H85").Value = 4787.125
$ws.Range("I85").Value = 1059.8
$ws.Range("K85").Value = 1059.8
$ws.Range("M85").Value = 188.2
$ws.Range("H93").Value = 2528365.5
$ws.Range("I93").Value = 2698.2727
$ws.Range("J93").Value = 5054033
$ws.Range("K93").Value = 2698.2727
$ws.Range("L93").Value = 5054033
$ws.Range("M93").Value = -1450.2727
$ws.Range("N93").Value = -5056529
$ws.Range("H132").Value = 3772.0625
$ws.Range("I132").Value = 2392.8076
$ws.Range("K132").Value = 7178.4228
$ws.Range("M132").Value = -4648.4228

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 27498.8
$ws.Range("J63").Value = 27498.8
$ws.Range("L63").Value = 27498.8
$ws.Range("N63").Value = -28746.8
$ws.Range("H66").Value = 27498.8
$ws.Range("J66").Value = 27498.8
$ws.Range("L66").Value = 82496.39999999999
$ws.Range("N66").Value = -88736.39999999999
$ws.Range("H81").Value = 2169.625
$ws.Range("I81").Value = 2268.2856
$ws.Range("K81").Value = 4536.5712
$ws.Range("M81").Value = -3475.5712
$ws.Range("H84").Value = 2169.625
$ws.Range("I84").Value = 2268.2856
$ws.Range("K84").Value = 22682.856
$ws.Range("M84").Value = -17378.856
$ws.Range("H126").Value = 4355
$ws.Range("J126").Value = 2729.7
$ws.Range("L126").Value = 8189.099999999999
$ws.Range("N126").Value = -13129.1
$ws.Range("H132").Value = 229780.34
$ws.Range("I132").Value = 2549.5405
$ws.Range("J132").Value = 1430857.4
$ws.Range("K132").Value = 7648.6215
$ws.Range("L132").Value = 4292572.199999999
$ws.Range("M132").Value = -5118.6215
$ws.Range("N132").Value = -4297632.199999999
